$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
#    (shared across Overview/zh-cn/de-de sheets; update every cell that
#    currently shows it so the underlying shared string is replaced
#    in place for all of them)
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("C2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("B3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("C3").Value = "Handed back: in sync with en-US"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("C3").Value = "Handed back: in sync with en-US"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("C3").Value = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# 2. zh-cn sheet: populate "Latest Target File" (F) / "Latest Handback
#    File" (G) columns + the handback datetime (H)
# ---------------------------------------------------------------------
$wsZhCn.Range("F2").Value = "a.md"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/3b222529bb42505978d9384de5d77ff74ef4b045/e2e/a.md", "", "", "a.md") | Out-Null

$wsZhCn.Range("G2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5f55b9e661a1bd9d52dbf4f116f8e3dfe0b28eb3/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf") | Out-Null

$wsZhCn.Range("F3").Value = "a.md"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/3b222529bb42505978d9384de5d77ff74ef4b045/e2e/a.md", "", "", "a.md") | Out-Null

$wsZhCn.Range("G3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5f55b9e661a1bd9d52dbf4f116f8e3dfe0b28eb3/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf") | Out-Null

# zh-cn handback datetime stays "not yet handed back" in this diff (H2/H3 unchanged)

# ---------------------------------------------------------------------
# 3. de-de sheet: populate "Latest Target File" (F) / "Latest Handback
#    File" (G) columns + the handback datetime (H), which IS updated
#    here (de-de has been handed back)
# ---------------------------------------------------------------------
$wsDeDe.Range("F2").Value = "a.md"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/3b222529bb42505978d9384de5d77ff74ef4b045/e2e/a.md", "", "", "a.md") | Out-Null

$wsDeDe.Range("G2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/692a81d6c75237ee0430288985e4dee96b6131d8/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf") | Out-Null

$wsDeDe.Range("F3").Value = "a.md"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/3b222529bb42505978d9384de5d77ff74ef4b045/e2e/a.md", "", "", "a.md") | Out-Null

$wsDeDe.Range("G3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/692a81d6c75237ee0430288985e4dee96b6131d8/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf") | Out-Null

$wsDeDe.Range("H2").Value = "2016-03-18 04:13:19"
$wsDeDe.Range("H3").Value = "2016-03-18 04:13:19"
